# proj4_report.pptx - slide 1 "PlaceHolder 2" (subtitle) placeholder fill-in.
# Replaces the bracketed placeholder text for name / GT email / GT username /
# GT ID with the real values, matching the author's commit ("pt 1 2 3 4").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find the subtitle placeholder shape by name rather than a hard-coded index.
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "PlaceHolder 2") {
        $shape = $candidate
    }
}

$tr = $shape.TextFrame.TextRange

# Paragraph 1: [name] -> Anirudh Arunkumar (also normalizes lang to en-US)
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "Anirudh Arunkumar"
$para1.LanguageID = "en-US"

# Paragraph 2: [GT email] -> aarunkumar8@gatech.edu
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "aarunkumar8@gatech.edu"

# Paragraph 3: [GT username] -> aarunkumar8
$para3 = $tr.Paragraphs(3, 1)
$para3.Text = "aarunkumar8"

# Paragraph 4: [GT ID] -> 903572206
$para4 = $tr.Paragraphs(4, 1)
$para4.Text = "903572206"
